$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 4617
    3  = 4685
    4  = 5060
    5  = 5060
    6  = 5060
    7  = 5087
    8  = 5087
    9  = 5148
    10 = 5148
    11 = 5148
    12 = 5148
    13 = 5148
    14 = 5153
    15 = 5153
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
